$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.314.24"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "2.303.11"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'316.06"
$ws.Range("E5").Value = "  +1.17%  "
$ws.Range("D6").Value = "'103.47"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("D7").Value = "'0.624"
$ws.Range("E7").Value = "  -0.50%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "'39.69"
$ws.Range("E10").Value = "  -1.57%  "
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "'8.38"
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("D13").Value = "'0.107"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("D14").Value = "'0.963"
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("E15").Value = "  -2.24%  "
$ws.Range("D16").Value = "2.650.26"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "2.294.99"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").Value = "42.420.39"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("E19").Value = "  -1.92%  "
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").Value = "'73.50"
$ws.Range("E21").Value = "  -1.60%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "'276.82"
$ws.Range("E22").Value = "  +7.34%  "
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").Value = "'3.53"
$ws.Range("E23").Value = "  +2.08%  "
$ws.Range("D24").Value = "'11.37"
$ws.Range("E24").Value = "  +22.02%  "
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D27").Value = "'10.85"
$ws.Range("E27").Value = "  -1.29%  "
$ws.Range("E28").Value = "  +3.32%  "
$ws.Range("D29").Value = "'22.75"
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("D30").Value = "'37.15"
$ws.Range("E30").Value = "  +3.89%  "
$ws.Range("D31").Value = "'165.43"
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("E32").Value = "  -2.36%  "
$ws.Range("D33").Value = "'5.88"
$ws.Range("E33").Value = "  +0.80%  "
$ws.Range("E34").Value = "  +5.19%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.118"
$ws.Range("E35").Value = "  -0.64%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "'2.60"
$ws.Range("E36").Value = "  -10.92%  "
$ws.Range("E37").Value = "  +3.69%  "
$ws.Range("D38").Value = "'4.58"
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("E39").Value = "  +2.22%  "
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("E41").Value = "  +3.23%  "
$ws.Range("D42").Value = "'69.95"
$ws.Range("E42").Value = "  -2.90%  "
$ws.Range("D43").Value = "'94.74"
$ws.Range("E43").Value = "  -3.66%  "
$ws.Range("E44").Value = "  -0.57%  "
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").Value = "'81.28"
$ws.Range("E46").Value = "  +8.57%  "
$ws.Range("D47").Value = "'12.08"
$ws.Range("E47").Value = "  -2.29%  "
$ws.Range("D48").Value = "'112.92"
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("E49").Value = "  -0.88%  "
$ws.Range("D50").Value = "'5.23"
$ws.Range("E50").Value = "  -1.88%  "
$ws.Range("D51").Value = "1.589.50"
$ws.Range("E51").Value = "  +1.20%  "
